# Generate Report for Handoff
# Updates status/dates for the 98f3bff0 and d548ec63 rows across the
# Overview, zh-cn and de-de sheets to reflect that a new handoff was
# generated, and records the "handback file is not latest" error detail
# for those two files on the localization sheets.

$wb = $excel.ActiveWorkbook

$msg98f3 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c62b879ff52abcf8abd16f8d3082bed5d7144ba4/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md."
$msgd548 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c62b879ff52abcf8abd16f8d3082bed5d7144ba4/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 -> 98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-28 12:26:19"

# Row 5 -> d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-28 12:26:19"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4 -> 98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-08-28 12:26:14"
$wsZhCn.Range("P4").Value = $msg98f3

# Row 5 -> d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-08-28 12:26:14"
$wsZhCn.Range("P5").Value = $msgd548

# Widen the Error Detail column to fit the new long text
$wsZhCn.Columns("P:P").ColumnWidth = 39.15

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 -> 98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-08-28 12:26:19"
$wsDeDe.Range("P4").Value = $msg98f3

# Row 5 -> d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-08-28 12:26:19"
$wsDeDe.Range("P5").Value = $msgd548

# Widen the Error Detail column to fit the new long text
$wsDeDe.Columns("P:P").ColumnWidth = 39.15
